# FA10_TestData_ManuallyAddAssets_21C.xlsx
# "Add files via upload" / "Anu - FA files uploaded"
#
# The re-uploaded workbook strips the hard-coded Oracle Fusion login
# (URL / username / password) that used to live in A2:C2 of the
# "Input_Value" sheet, together with the hyperlink that decorated the
# URL cell. Everything else on that sheet (headers, remaining sample
# values) is left exactly as it was.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Input_Value")
$ws.Activate()

# Remove the hyperlink that lived on A2 (https://edrx.fa.us2.oraclecloud.com/)
foreach ($hl in $ws.Hyperlinks) {
    $hl.Delete()
}

# Clear the stored credentials: UserName (A2), Password (B2) and
# SelectBook (C2) - values only, cell styles are left untouched.
$ws.Range("A2").Value = ""
$ws.Range("B2").Value = ""
$ws.Range("C2").Value = ""

# Restore the on-screen selection to A2:C2 (the cells that were just
# cleared), matching the saved sheet view.
[void]$ws.Range("A2:C2").Select()
